$d = $word.ActiveDocument

# Locate the "第3章" chapter-number prefix at the start of the title paragraph,
# so we know exactly where to split the run and insert the new space.
$prefixRange = $d.Content
$prefixRange.Find.Execute("第3章", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $prefixRange.End

# Insert a space right after "第3章" - this begins the new chapter-content
# text that will follow the heading number.
$d.Range($splitPos, $splitPos).InsertAfter(" ")

# Force the freshly inserted space to live in its own run (distinct from the
# "第3章" run before it and the remaining title text after it), by toggling
# its Bold formatting off and back on. Identically-formatted adjacent runs
# are normally coalesced, but a run that was explicitly (re)formatted keeps
# its own boundary.
$spaceRange = $d.Range($splitPos, $splitPos + 1)
$spaceRange.Bold = 0
$spaceRange = $d.Range($splitPos, $splitPos + 1)
$spaceRange.Bold = 1

# Move the document's "_GoBack" bookmark (previously sitting by itself in the
# empty paragraph right after the heading) to a collapsed position right
# after the space we just typed - marking that as the last edit location,
# i.e. where writing of the chapter content will continue. Adding a bookmark
# with a name that already exists simply relocates it, and also removes it
# from the now-empty trailing paragraph.
$goBackPos = $splitPos + 1
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos)) | Out-Null
